$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 1 values (F1, G1, H1 change; I1 and J1 are removed entirely)
$ws.Range("F1").Value = 3
$ws.Range("G1").Value = 1
$ws.Range("H1").Value = 73
$ws.Range("I1").ClearContents()
$ws.Range("J1").ClearContents()

# Add a new value in row 7
$ws.Range("B7").Value = 22
